$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.270036666666666
$ws.Range("H2").Value = 9.81011
$ws.Range("I2").Value = 0.359406393324744
$ws.Range("J2").Value = 0.3594063933247441
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1030763333333333
$ws.Range("N2").Value = 0.309229
$ws.Range("O2").Value = 0.01126512502660735
$ws.Range("P2").Value = 0.01126512502660735
$ws.Range("Q2").Value = 0.3370633894655555
$ws.Range("R2").Value = 3.03357050519
$ws.Range("S2").Value = 0.004048757956165259
$ws.Range("T2").Value = 0.004048757956165261
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.270036666666666
$ws.Range("H3").Value = 9.81011
$ws.Range("I3").Value = 0.359406393324744
$ws.Range("J3").Value = 0.3594063933247441
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.046962666666667
$ws.Range("N3").Value = 27.140888
$ws.Range("O3").Value = 0.9887348749733926
$ws.Range("P3").Value = 0.9887348749733927
$ws.Range("Q3").Value = 29.58389964196445
$ws.Range("R3").Value = 266.25509677768
$ws.Range("S3").Value = 0.3553576353685787
$ws.Range("T3").Value = 0.3553576353685788
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.828401
$ws.Range("H4").Value = 17.485203
$ws.Range("I4").Value = 0.6405936066752559
$ws.Range("J4").Value = 0.640593606675256
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1030763333333333
$ws.Range("N4").Value = 0.309229
$ws.Range("O4").Value = 0.01126512502660735
$ws.Range("P4").Value = 0.01126512502660735
$ws.Range("Q4").Value = 0.6007702042763333
$ws.Range("R4").Value = 5.406931838487
$ws.Range("S4").Value = 0.007216367070442092
$ws.Range("T4").Value = 0.007216367070442094
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.828401
$ws.Range("H5").Value = 17.485203
$ws.Range("I5").Value = 0.6405936066752559
$ws.Range("J5").Value = 0.640593606675256
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.046962666666667
$ws.Range("N5").Value = 27.140888
$ws.Range("O5").Value = 0.9887348749733926
$ws.Range("P5").Value = 0.9887348749733927
$ws.Range("Q5").Value = 52.72932625336267
$ws.Range("R5").Value = 474.563936280264
$ws.Range("S5").Value = 0.6333772396048137
$ws.Range("T5").Value = 0.633377239604814
